$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update existing rows with corrected figures
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 187889
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 20877

$ws.Range("D6").Value = 265513
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 17701

$ws.Range("D10").Value = 205890
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = 12868

# ---------------------------------------------------------------------------
# 2) Append newly extracted rows (15-20) below the existing table data,
#    copying the formatting from the last existing row (14) so number
#    formats / styles match the rest of the table.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 15; Date = 43538; Prod = 456763; Res = 21; Ippr = 21751 },
    @{ Row = 16; Date = 43539; Prod = 434511; Res = 20; Ippr = 21726 },
    @{ Row = 17; Date = 43540; Prod = 124044; Res = 9;  Ippr = 13783 },
    @{ Row = 18; Date = 43541; Prod = 19498;  Res = 2;  Ippr = 9749 },
    @{ Row = 19; Date = 43542; Prod = 210018; Res = 14; Ippr = 15001 },
    @{ Row = 20; Date = 43543; Prod = 153867; Res = 10; Ippr = 15387 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy formatting (number formats / styles) from the row above so the
    # new row matches the rest of the table.
    $ws.Range("A" + ($row - 1) + ":F" + ($row - 1)).Copy() | Out-Null
    $ws.Range("A" + $row + ":F" + $row).PasteSpecial(-4122) | Out-Null

    $ws.Range("A" + $row).Value = 2019
    $ws.Range("B" + $row).Value = "marzo"
    $ws.Range("C" + $row).Value = $r.Date
    $ws.Range("D" + $row).Value = $r.Prod
    $ws.Range("E" + $row).Value = $r.Res
    $ws.Range("F" + $row).Value = $r.Ippr
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Grow the Table1 ListObject so the new rows are recognised as part of
#    the table (ref A1:F14 -> A1:F20).
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:F20"))
